$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (data starts on row 2, header on row 1).
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 284) {
    $lastRow = 284
}

# Column C holds the "Förändrad" (changed) date for every record.
# Bump every row's date-serial from 46061 (2026-02-08) to 46062 (2026-02-09).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
